$p = $ppt.ActivePresentation

# Slide 6 ("Title 1" shape): "Problem Statement" -> "Existing System"
# The original paragraph also carries a stray <a:endParaRPr> run-properties
# element (left over from a prior edit) which the real authoring app
# removes once the text is reset/retyped. Deleting the existing text range
# first (which also clears the paragraph's trailing run-properties) and
# then typing the new text reproduces that same cleaned-up shape.
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(1)
$shp6.TextFrame.TextRange.Delete()
$shp6.TextFrame.TextRange.Text = "Existing System"

# Slide 7 ("Title 1" shape): "Solution Overview" -> "Proposed System"
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(1)
$shp7.TextFrame.TextRange.Text = "Proposed System"
